$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, pushing the existing rows 37-38 down to 38-39.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the latest weekly price data.
$ws.Cells.Item(37, 1).Value = 11
$ws.Cells.Item(37, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(37, 3).Value = "Bíobío"
$ws.Cells.Item(37, 4).Value = 44706
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
$ws.Cells.Item(37, 5).Value = 8
$ws.Cells.Item(37, 6).Value = 100112013
$ws.Cells.Item(37, 7).Value = "Alcachofa"
$ws.Cells.Item(37, 8).Value = "Española"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 50
$ws.Cells.Item(37, 11).Value = 23000
$ws.Cells.Item(37, 12).Value = 24000
$ws.Cells.Item(37, 13).Value = 23400
$ws.Cells.Item(37, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(37, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(37, 16).Value = 780
$ws.Cells.Item(37, 17).Value = 30
$ws.Cells.Item(37, 18).Value = "Hortaliza"
